# "The Last Update 15-03-2024"
# Refresh the Super Lig standings table (teams C-G stat columns + re-ranked team names)
# with the latest odds/percentages, mirroring the new shared-string order from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the plain (unstyled) look of the data cells so it can be restored at the end.
$baseStyle = $ws.Range("A2").Style

# Temporarily mark the data block as Text so every value we write below is stored
# as a literal string (matching the workbook`'s original t="s" shared-string cells)
# instead of being auto-converted into numbers/percentages.
$dataRange = $ws.Range("B2:G21")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "Galatasaray"
$ws.Range("C2").Value = "1.9"
$ws.Range("D2").Value = "7.4"
$ws.Range("E2").Value = "79%"
$ws.Range("F2").Value = "57%"
$ws.Range("G2").Value = "2.61"

$ws.Range("B3").Value = "Fenerbahçe"
$ws.Range("C3").Value = "2.5"
$ws.Range("D3").Value = "6.3"
$ws.Range("E3").Value = "79%"
$ws.Range("F3").Value = "61%"
$ws.Range("G3").Value = "3.36"

$ws.Range("B4").Value = "Trabzonspor"
$ws.Range("C4").Value = "1.7"
$ws.Range("D4").Value = "4.6"
$ws.Range("E4").Value = "69%"
$ws.Range("F4").Value = "46%"
$ws.Range("G4").Value = "2.81"

$ws.Range("B5").Value = "Beşiktaş"
$ws.Range("C5").Value = "2.3"
$ws.Range("D5").Value = "5.7"
$ws.Range("E5").Value = "79%"
$ws.Range("F5").Value = "46%"
$ws.Range("G5").Value = "2.54"

$ws.Range("B6").Value = "Kasımpaşa"
$ws.Range("C6").Value = "1.9"
$ws.Range("D6").Value = "5.0"
$ws.Range("E6").Value = "86%"
$ws.Range("F6").Value = "71%"
$ws.Range("G6").Value = "3.50"

$ws.Range("B7").Value = "Başakşehir"
$ws.Range("C7").Value = "2.1"
$ws.Range("D7").Value = "4.0"
$ws.Range("E7").Value = "68%"
$ws.Range("F7").Value = "36%"
$ws.Range("G7").Value = "2.46"

$ws.Range("B8").Value = "Rizespor"
$ws.Range("C8").Value = "2.4"
$ws.Range("D8").Value = "4.8"
$ws.Range("E8").Value = "66%"
$ws.Range("F8").Value = "51%"
$ws.Range("G8").Value = "2.54"

$ws.Range("B9").Value = "Sivasspor"
$ws.Range("C9").Value = "2.0"
$ws.Range("D9").Value = "3.6"
$ws.Range("E9").Value = "71%"
$ws.Range("F9").Value = "42%"
$ws.Range("G9").Value = "2.54"

$ws.Range("B10").Value = "Antalyaspor"
$ws.Range("C10").Value = "2.1"
$ws.Range("D10").Value = "5.6"
$ws.Range("E10").Value = "74%"
$ws.Range("F10").Value = "37%"
$ws.Range("G10").Value = "2.33"

$ws.Range("B11").Value = "Kayserispor"
$ws.Range("C11").Value = "2.1"
$ws.Range("D11").Value = "4.5"
$ws.Range("E11").Value = "79%"
$ws.Range("F11").Value = "46%"
$ws.Range("G11").Value = "2.71"

$ws.Range("B12").Value = "Samsunspor"
$ws.Range("C12").Value = "2.1"
$ws.Range("D12").Value = "4.7"
$ws.Range("E12").Value = "82%"
$ws.Range("F12").Value = "39%"
$ws.Range("G12").Value = "2.50"

$ws.Range("B13").Value = "Adana Demirspor"
$ws.Range("C13").Value = "2.2"
$ws.Range("D13").Value = "5.0"
$ws.Range("E13").Value = "74%"
$ws.Range("F13").Value = "56%"
$ws.Range("G13").Value = "2.78"

$ws.Range("B14").Value = "Konyaspor"
$ws.Range("C14").Value = "2.0"
$ws.Range("D14").Value = "4.1"
$ws.Range("E14").Value = "82%"
$ws.Range("F14").Value = "32%"
$ws.Range("G14").Value = "2.39"

$ws.Range("B15").Value = "Ankaragücü"
$ws.Range("C15").Value = "1.6"
$ws.Range("D15").Value = "3.5"
$ws.Range("E15").Value = "79%"
$ws.Range("F15").Value = "40%"
$ws.Range("G15").Value = "2.46"

$ws.Range("B16").Value = "Alanyaspor"
$ws.Range("C16").Value = "2.6"
$ws.Range("D16").Value = "5.0"
$ws.Range("E16").Value = "74%"
$ws.Range("F16").Value = "52%"
$ws.Range("G16").Value = "2.67"

$ws.Range("B17").Value = "Karagümrük"
$ws.Range("C17").Value = "2.4"
$ws.Range("D17").Value = "4.9"
$ws.Range("E17").Value = "78%"
$ws.Range("F17").Value = "41%"
$ws.Range("G17").Value = "2.33"

$ws.Range("B18").Value = "Hatayspor"
$ws.Range("C18").Value = "2.8"
$ws.Range("D18").Value = "4.6"
$ws.Range("E18").Value = "71%"
$ws.Range("F18").Value = "46%"
$ws.Range("G18").Value = "2.54"

$ws.Range("B19").Value = "Gaziantep"
$ws.Range("C19").Value = "2.4"
$ws.Range("D19").Value = "4.0"
$ws.Range("E19").Value = "85%"
$ws.Range("F19").Value = "49%"
$ws.Range("G19").Value = "2.63"

$ws.Range("B20").Value = "Pendikspor"
$ws.Range("C20").Value = "2.0"
$ws.Range("D20").Value = "4.6"
$ws.Range("E20").Value = "93%"
$ws.Range("F20").Value = "50%"
$ws.Range("G20").Value = "3.18"

$ws.Range("B21").Value = "İstanbulspor"
$ws.Range("C21").Value = "2.4"
$ws.Range("D21").Value = "4.1"
$ws.Range("E21").Value = "78%"
$ws.Range("F21").Value = "45%"
$ws.Range("G21").Value = "2.59"

# Put the original (default) style back now that the text values are committed.
$dataRange.Style = $baseStyle
